$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.139.93'
$ws.Range('E2').Value = '  +0.46%  '

# Row 3
$ws.Range('D3').Value = '1.920.10'

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.08%  '

# Row 6
$ws.Range('E6').Value = '  +0.05%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5070'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '

# Row 8
$ws.Range('E8').Value = '  +3.30%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08357'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.15%  '

# Row 10
$ws.Range('E10').Value = '  +1.54%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.31'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.36%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.00'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.52%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.434'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.68%  '

# Row 14
$ws.Range('D14').Value = '1.916.47'
$ws.Range('E14').Value = '  +2.67%  '

# Row 15
$ws.Range('E15').Value = '  +1.18%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.01%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.57'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.57%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001097'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.01%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06516'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.32%  '

# Row 20
$ws.Range('E20').Value = '  +3.36%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.952'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.48%  '

# Row 23
$ws.Range('D23').Value = '30.153.18'
$ws.Range('E23').Value = '  +0.57%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.37%  '

# Row 25
$ws.Range('E25').Value = '  +2.18%  '

# Row 26
$ws.Range('D26').Value = '2.135.59'
$ws.Range('E26').Value = '  +2.50%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.41%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.06%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.265'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.33%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '129.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.57%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.131'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.45%  '

# Row 32
$ws.Range('E32').Value = '  +0.97%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.959'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.93%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.797'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.80%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02454'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.33%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.330'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.01%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06441'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.49%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2151'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.30%  '

# Row 39
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.199'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.20%  '

# Row 40
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6512'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.22%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.592'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.24%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.45'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.95%  '

# Row 43
$ws.Range('E43').Value = '  +1.41%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.15%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.188'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +9.66%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6057'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.65%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.624'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.18%  '

# Row 48
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.13%  '

# Row 49
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.211'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.13%  '

# Row 50
$ws.Range('E50').Value = '  +1.82%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.00%  '
